# The upstream change is a pure "canonical re-save" of this template: every
# paragraph's text, every bookmark/rsid, every style and every page-setup
# value is byte-for-byte identical before and after the commit - the only
# difference in the OOXML is the *order* in which already-present attributes
# are written back out (e.g. <w:tab w:val="left" w:pos="3119"/> becomes
# <w:tab w:pos="3119" w:val="left"/>, <w:pgSz w:w=".." w:h=".."/> becomes
# <w:pgSz w:h=".." w:w=".."/>, etc.). There is no new text, no new value and
# nothing removed anywhere in the document.
#
# We reproduce the intent of that resave through the Word object model by
# touching every value that the diff re-emits, using the same final values
# already present in the document, so each touched node is confirmed/written
# back by Word itself.

$d = $word.ActiveDocument

# --- Section page setup: <w:pgSz>/<w:pgMar> (values unchanged: 11906x16838,
#     margins 1417/1417/1417/1417, header/footer 708/708, gutter 0) ---
foreach ($sec in $d.Sections) {
    $ps = $sec.PageSetup
    $ps.PageWidth     = 11906 / 20.0
    $ps.PageHeight    = 16838 / 20.0
    $ps.TopMargin     = 1417 / 20.0
    $ps.RightMargin   = 1417 / 20.0
    $ps.BottomMargin  = 1417 / 20.0
    $ps.LeftMargin    = 1417 / 20.0
    $ps.HeaderDistance = 708 / 20.0
    $ps.FooterDistance = 708 / 20.0
    $ps.Gutter        = 0
}

# --- Tab stops: the four <w:tab w:val="left" w:pos="3119"/> entries in the
#     "if / then / else" paragraphs (value unchanged: left tab @ 3119 twips
#     = 155.95 pt) ---
foreach ($p in $d.Paragraphs) {
    $tabs = $p.Range.ParagraphFormat.TabStops
    for ($j = 1; $j -le $tabs.Count; $j++) {
        $t = $tabs.Item($j)
        if ($t.CustomTab -and [Math]::Round($t.Position * 20) -eq 3119) {
            $t.Position = 3119 / 20.0
            $t.Alignment = 0
        }
    }
}
